$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) TestRecord sheet: extend the "some test text..." note by three more z's,
#    and bump the row-10 date/amount.
# ---------------------------------------------------------------------------
$testRecord = $wb.Worksheets.Item("TestRecord")
$testRecord.Range("E10").Value = "some test textzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"
$testRecord.Range("A10").Value = 43233
$testRecord.Range("B10").Value = 85.14

# ---------------------------------------------------------------------------
# 2) Expected Out sheet: bump a contributing amount; the SUM(B2:B295) formula
#    in B1 recalculates automatically.
# ---------------------------------------------------------------------------
$expectedOut = $wb.Worksheets.Item("Expected Out")
$expectedOut.Range("B9").Value = 1342.16

# ---------------------------------------------------------------------------
# 3) CredCard sheet: bring in reconciliation columns F:J (cheque number,
#    reconciled amount, notes/info) mirroring the layout already used on the
#    Bank sheet, across the existing data rows.
# ---------------------------------------------------------------------------
$credCard = $wb.Worksheets.Item("CredCard")

$credCard.Range("F1").Value = "Cheque num"
$credCard.Range("F1").NumberFormat = "0"
$credCard.Range("G1").Value = "Reconciled Amt"
$credCard.Range("G1").NumberFormat = '"£"#,##0.00_);[Red]\("£"#,##0.00\)'

function Set-CCRow($row, $f, $g, $h, $includeI, $j) {
    $credCard.Range("F$row").Value = $f
    $credCard.Range("F$row").NumberFormat = "0"
    $credCard.Range("G$row").Value = $g
    $credCard.Range("G$row").NumberFormat = '"£"#,##0.00_);[Red]\("£"#,##0.00\)'
    $credCard.Range("H$row").Value = $h
    $credCard.Range("H$row").NumberFormat = "General"
    if ($includeI) {
        $credCard.Range("I$row").NumberFormat = "General"
    }
    $credCard.Range("J$row").Value = $j
    $credCard.Range("J$row").NumberFormat = "General"
}

Set-CCRow 3  12345 4567.8900000000003 "notes"  $true  "info"
Set-CCRow 4  12345 4567.8900000000003 "notes1" $true  "info1"
Set-CCRow 6  22345 5567.89             "notes2" $false "info2"
Set-CCRow 7  42345 7567.89             "notes4" $false "info4"
Set-CCRow 8  42345 7567.89             "notes4" $false "info4"
Set-CCRow 10 42345 7567.89             "notes4" $false "info4"
Set-CCRow 11 42345 7567.89             "notes4" $false "info4"

$credCard.Columns.Item(7).ColumnWidth = 14.140625

# ---------------------------------------------------------------------------
# 4) View-state: selections on Bank / CredCard, and CredCard becomes the
#    active sheet/tab.
# ---------------------------------------------------------------------------
$bank = $wb.Worksheets.Item("Bank")
$bank.Activate() | Out-Null
$bank.Range("F1:J7").Select() | Out-Null

$credCard.Activate() | Out-Null
$credCard.Range("F11:J11").Select() | Out-Null
